$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C (old C -> D), keeping A/B untouched.
$ws.Columns("C:C").Insert()

# Insert a new column at E for the "Type" classification.
$ws.Columns("E:E").Insert()

# Header row
$ws.Range("C1").Value = "Intercept"
$ws.Range("E1").Value = "Type"

# New numeric "Intercept" column values (only some rows have a value)
$ws.Range("C3").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("C5").Value = 5
$ws.Range("C6").Value = 5
$ws.Range("C7").Value = 3
$ws.Range("C8").Value = 0

# New "Type" column values
$ws.Range("E2").Value = "Bool"
$ws.Range("E3").Value = "Int"
$ws.Range("E4").Value = "Int"
$ws.Range("E5").Value = "Int"
$ws.Range("E6").Value = "Int"
$ws.Range("E7").Value = "Int"
$ws.Range("E8").Value = "Int"
$ws.Range("E9").Value = "Bool"
$ws.Range("E10").Value = "Bool"
$ws.Range("E11").Value = "Bool"
$ws.Range("E12").Value = "Bool"

# Update the selection to match the saved workbook state
$null = $ws.Range("C7").Select()
